$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 43.34730933333333
$ws.Range("H2").Value = 130.041928
$ws.Range("I2").Value = 0.04273139820300816
$ws.Range("J2").Value = 0.04273139820300816
$ws.Range("M2").Value = 18.95273633333333
$ws.Range("N2").Value = 56.858209
$ws.Range("O2").Value = 0.08721078561875104
$ws.Range("P2").Value = 0.08721078561875105
$ws.Range("Q2").Value = 821.5501245541055
$ws.Range("R2").Value = 7393.951120986951
$ws.Range("S2").Value = 0.003726638807872028
$ws.Range("T2").Value = 0.003726638807872028
$ws.Range("G3").Value = 43.34730933333333
$ws.Range("H3").Value = 130.041928
$ws.Range("I3").Value = 0.04273139820300816
$ws.Range("J3").Value = 0.04273139820300816
$ws.Range("O3").Value = 0.04852204497892696
$ws.Range("P3").Value = 0.04852204497892696
$ws.Range("Q3").Value = 457.0913082966933
$ws.Range("R3").Value = 4113.821774670239
$ws.Range("S3").Value = 0.002073414825618801
$ws.Range("T3").Value = 0.0020734148256188
$ws.Range("G4").Value = 43.34730933333333
$ws.Range("H4").Value = 130.041928
$ws.Range("I4").Value = 0.04273139820300816
$ws.Range("J4").Value = 0.04273139820300816
$ws.Range("M4").Value = 101.4555613333333
$ws.Range("N4").Value = 304.366684
$ws.Range("O4").Value = 0.4668465309523581
$ws.Range("P4").Value = 0.4668465309523581
$ws.Range("Q4").Value = 4397.825600702972
$ws.Range("R4").Value = 39580.43040632675
$ws.Range("S4").Value = 0.01994900501381819
$ws.Range("T4").Value = 0.01994900501381819
$ws.Range("G5").Value = 43.34730933333333
$ws.Range("H5").Value = 130.041928
$ws.Range("I5").Value = 0.04273139820300816
$ws.Range("J5").Value = 0.04273139820300816
$ws.Range("M5").Value = 2.410466333333333
$ws.Range("N5").Value = 7.231399000000001
$ws.Range("O5").Value = 0.01109173150200089
$ws.Range("P5").Value = 0.01109173150200089
$ws.Range("Q5").Value = 104.4872297885858
$ws.Range("R5").Value = 940.385068097272
$ws.Range("S5").Value = 0.00047396519557285
$ws.Range("T5").Value = 0.00047396519557285
$ws.Range("G6").Value = 43.34730933333333
$ws.Range("H6").Value = 130.041928
$ws.Range("I6").Value = 0.04273139820300816
$ws.Range("J6").Value = 0.04273139820300816
$ws.Range("M6").Value = 83.95738966666666
$ws.Range("N6").Value = 251.872169
$ws.Range("O6").Value = 0.386328906947963
$ws.Range("P6").Value = 0.386328906947963
$ws.Range("Q6").Value = 3639.326940700203
$ws.Range("R6").Value = 32753.94246630183
$ws.Range("S6").Value = 0.0165083743601263
$ws.Range("T6").Value = 0.01650837436012629
$ws.Range("I7").Value = 0.0889365509391893
$ws.Range("J7").Value = 0.08893655093918929
$ws.Range("M7").Value = 18.95273633333333
$ws.Range("N7").Value = 56.858209
$ws.Range("O7").Value = 0.08721078561875104
$ws.Range("P7").Value = 0.08721078561875105
$ws.Range("Q7").Value = 1709.886349947705
$ws.Range("R7").Value = 15388.97714952934
$ws.Range("S7").Value = 0.007756226477628769
$ws.Range("T7").Value = 0.007756226477628769
$ws.Range("I8").Value = 0.0889365509391893
$ws.Range("J8").Value = 0.08893655093918929
$ws.Range("O8").Value = 0.04852204497892696
$ws.Range("P8").Value = 0.04852204497892696
$ws.Range("S8").Value = 0.004315383324941972
$ws.Range("T8").Value = 0.004315383324941972
$ws.Range("I9").Value = 0.0889365509391893
$ws.Range("J9").Value = 0.08893655093918929
$ws.Range("M9").Value = 101.4555613333333
$ws.Range("N9").Value = 304.366684
$ws.Range("O9").Value = 0.4668465309523581
$ws.Range("P9").Value = 0.4668465309523581
$ws.Range("Q9").Value = 9153.162709547296
$ws.Range("R9").Value = 82378.46438592568
$ws.Range("S9").Value = 0.04151972028082821
$ws.Range("T9").Value = 0.0415197202808282
$ws.Range("I10").Value = 0.0889365509391893
$ws.Range("J10").Value = 0.08893655093918929
$ws.Range("M10").Value = 2.410466333333333
$ws.Range("N10").Value = 7.231399000000001
$ws.Range("O10").Value = 0.01109173150200089
$ws.Range("P10").Value = 0.01109173150200089
$ws.Range("Q10").Value = 217.4685178902748
$ws.Range("R10").Value = 1957.216661012473
$ws.Range("S10").Value = 0.0009864603437315133
$ws.Range("T10").Value = 0.000986460343731513
$ws.Range("I11").Value = 0.0889365509391893
$ws.Range("J11").Value = 0.08893655093918929
$ws.Range("M11").Value = 83.95738966666666
$ws.Range("N11").Value = 251.872169
$ws.Range("O11").Value = 0.386328906947963
$ws.Range("P11").Value = 0.386328906947963
$ws.Range("Q11").Value = 7574.504918099361
$ws.Range("R11").Value = 68170.54426289425
$ws.Range("S11").Value = 0.03435876051205883
$ws.Range("T11").Value = 0.03435876051205883
$ws.Range("G12").Value = 394.701121
$ws.Range("H12").Value = 1184.103363
$ws.Range("I12").Value = 0.3890929110023202
$ws.Range("J12").Value = 0.3890929110023201
$ws.Range("M12").Value = 18.95273633333333
$ws.Range("N12").Value = 56.858209
$ws.Range("O12").Value = 0.08721078561875104
$ws.Range("P12").Value = 0.08721078561875105
$ws.Range("Q12").Value = 7480.666276784095
$ws.Range("R12").Value = 67325.99649105685
$ws.Range("S12").Value = 0.03393309844719912
$ws.Range("T12").Value = 0.03393309844719912
$ws.Range("G13").Value = 394.701121
$ws.Range("H13").Value = 1184.103363
$ws.Range("I13").Value = 0.3890929110023202
$ws.Range("J13").Value = 0.3890929110023201
$ws.Range("O13").Value = 0.04852204497892696
$ws.Range("P13").Value = 0.04852204497892696
$ws.Range("Q13").Value = 4162.06806278806
$ws.Range("R13").Value = 37458.61256509254
$ws.Range("S13").Value = 0.01887958372863621
$ws.Range("T13").Value = 0.0188795837286362
$ws.Range("G14").Value = 394.701121
$ws.Range("H14").Value = 1184.103363
$ws.Range("I14").Value = 0.3890929110023202
$ws.Range("J14").Value = 0.3890929110023201
$ws.Range("M14").Value = 101.4555613333333
$ws.Range("N14").Value = 304.366684
$ws.Range("O14").Value = 0.4668465309523581
$ws.Range("P14").Value = 0.4668465309523581
$ws.Range("Q14").Value = 40044.62378995092
$ws.Range("R14").Value = 360401.6141095583
$ws.Range("S14").Value = 0.1816466757195878
$ws.Range("T14").Value = 0.1816466757195877
$ws.Range("G15").Value = 394.701121
$ws.Range("H15").Value = 1184.103363
$ws.Range("I15").Value = 0.3890929110023202
$ws.Range("J15").Value = 0.3890929110023201
$ws.Range("M15").Value = 2.410466333333333
$ws.Range("N15").Value = 7.231399000000001
$ws.Range("O15").Value = 0.01109173150200089
$ws.Range("P15").Value = 0.01109173150200089
$ws.Range("Q15").Value = 951.4137638994264
$ws.Range("R15").Value = 8562.723875094838
$ws.Range("S15").Value = 0.004315714098169665
$ws.Range("T15").Value = 0.004315714098169664
$ws.Range("G16").Value = 394.701121
$ws.Range("H16").Value = 1184.103363
$ws.Range("I16").Value = 0.3890929110023202
$ws.Range("J16").Value = 0.3890929110023201
$ws.Range("M16").Value = 83.95738966666666
$ws.Range("N16").Value = 251.872169
$ws.Range("O16").Value = 0.386328906947963
$ws.Range("P16").Value = 0.386328906947963
$ws.Range("Q16").Value = 33138.07581766714
$ws.Range("R16").Value = 298242.6823590043
$ws.Range("S16").Value = 0.1503178390087274
$ws.Range("T16").Value = 0.1503178390087274
$ws.Range("G17").Value = 7.804371333333333
$ws.Range("H17").Value = 23.413114
$ws.Range("I17").Value = 0.007693480963358413
$ws.Range("J17").Value = 0.007693480963358412
$ws.Range("M17").Value = 18.95273633333333
$ws.Range("N17").Value = 56.858209
$ws.Range("O17").Value = 0.08721078561875104
$ws.Range("P17").Value = 0.08721078561875105
$ws.Range("Q17").Value = 147.9141921280917
$ws.Range("R17").Value = 1331.227729152826
$ws.Range("S17").Value = 0.0006709545189573928
$ws.Range("T17").Value = 0.0006709545189573928
$ws.Range("G18").Value = 7.804371333333333
$ws.Range("H18").Value = 23.413114
$ws.Range("I18").Value = 0.007693480963358413
$ws.Range("J18").Value = 0.007693480963358412
$ws.Range("O18").Value = 0.04852204497892696
$ws.Range("P18").Value = 0.04852204497892696
$ws.Range("Q18").Value = 82.29600309801333
$ws.Range("R18").Value = 740.66402788212
$ws.Range("S18").Value = 0.0003733034293485952
$ws.Range("T18").Value = 0.0003733034293485952
$ws.Range("G19").Value = 7.804371333333333
$ws.Range("H19").Value = 23.413114
$ws.Range("I19").Value = 0.007693480963358413
$ws.Range("J19").Value = 0.007693480963358412
$ws.Range("M19").Value = 101.4555613333333
$ws.Range("N19").Value = 304.366684
$ws.Range("O19").Value = 0.4668465309523581
$ws.Range("P19").Value = 0.4668465309523581
$ws.Range("Q19").Value = 791.7968744771084
$ws.Range("R19").Value = 7126.171870293976
$ws.Range("S19").Value = 0.003591674898691881
$ws.Range("T19").Value = 0.003591674898691881
$ws.Range("G20").Value = 7.804371333333333
$ws.Range("H20").Value = 23.413114
$ws.Range("I20").Value = 0.007693480963358413
$ws.Range("J20").Value = 0.007693480963358412
$ws.Range("M20").Value = 2.410466333333333
$ws.Range("N20").Value = 7.231399000000001
$ws.Range("O20").Value = 0.01109173150200089
$ws.Range("P20").Value = 0.01109173150200089
$ws.Range("Q20").Value = 18.81217435183178
$ws.Range("R20").Value = 169.309569166486
$ws.Range("S20").Value = 0.0000853340251613267
$ws.Range("T20").Value = 0.00008533402516132668
$ws.Range("G21").Value = 7.804371333333333
$ws.Range("H21").Value = 23.413114
$ws.Range("I21").Value = 0.007693480963358413
$ws.Range("J21").Value = 0.007693480963358412
$ws.Range("M21").Value = 83.95738966666666
$ws.Range("N21").Value = 251.872169
$ws.Range("O21").Value = 0.386328906947963
$ws.Range("P21").Value = 0.386328906947963
$ws.Range("Q21").Value = 655.2346451360295
$ws.Range("R21").Value = 5897.111806224266
$ws.Range("S21").Value = 0.002972214091199217
$ws.Range("T21").Value = 0.002972214091199217
$ws.Range("G22").Value = 478.3423056666667
$ws.Range("H22").Value = 1435.026917
$ws.Range("I22").Value = 0.4715456588921241
$ws.Range("J22").Value = 0.471545658892124
$ws.Range("M22").Value = 18.95273633333333
$ws.Range("N22").Value = 56.858209
$ws.Range("O22").Value = 0.08721078561875104
$ws.Range("P22").Value = 0.08721078561875105
$ws.Range("Q22").Value = 9065.895596379072
$ws.Range("R22").Value = 81593.06036741166
$ws.Range("S22").Value = 0.04112386736709374
$ws.Range("T22").Value = 0.04112386736709373
$ws.Range("G23").Value = 478.3423056666667
$ws.Range("H23").Value = 1435.026917
$ws.Range("I23").Value = 0.4715456588921241
$ws.Range("J23").Value = 0.471545658892124
$ws.Range("O23").Value = 0.04852204497892696
$ws.Range("P23").Value = 0.04852204497892696
$ws.Range("Q23").Value = 5044.052645332207
$ws.Range("R23").Value = 45396.47380798987
$ws.Range("S23").Value = 0.0228803596703814
$ws.Range("T23").Value = 0.02288035967038139
$ws.Range("G24").Value = 478.3423056666667
$ws.Range("H24").Value = 1435.026917
$ws.Range("I24").Value = 0.4715456588921241
$ws.Range("J24").Value = 0.471545658892124
$ws.Range("M24").Value = 101.4555613333333
$ws.Range("N24").Value = 304.366684
$ws.Range("O24").Value = 0.4668465309523581
$ws.Range("P24").Value = 0.4668465309523581
$ws.Range("Q24").Value = 48530.48713089259
$ws.Range("R24").Value = 436774.3841780333
$ws.Range("S24").Value = 0.2201394550394321
$ws.Range("T24").Value = 0.2201394550394321
$ws.Range("G25").Value = 478.3423056666667
$ws.Range("H25").Value = 1435.026917
$ws.Range("I25").Value = 0.4715456588921241
$ws.Range("J25").Value = 0.471545658892124
$ws.Range("M25").Value = 2.410466333333333
$ws.Range("N25").Value = 7.231399000000001
$ws.Range("O25").Value = 0.01109173150200089
$ws.Range("P25").Value = 0.01109173150200089
$ws.Range("Q25").Value = 1153.028023618543
$ws.Range("R25").Value = 10377.25221256689
$ws.Range("S25").Value = 0.005230257839365541
$ws.Range("T25").Value = 0.005230257839365539
$ws.Range("G26").Value = 478.3423056666667
$ws.Range("H26").Value = 1435.026917
$ws.Range("I26").Value = 0.4715456588921241
$ws.Range("J26").Value = 0.471545658892124
$ws.Range("M26").Value = 83.95738966666666
$ws.Range("N26").Value = 251.872169
$ws.Range("O26").Value = 0.386328906947963
$ws.Range("P26").Value = 0.386328906947963
$ws.Range("Q26").Value = 40160.37135090811
$ws.Range("R26").Value = 361443.342158173
$ws.Range("S26").Value = 0.1821717189758513
$ws.Range("T26").Value = 0.1821717189758513
